# Solve Leetcode - 1448. Count Good Nodes in Binary Tree - DFS
# Adds a new row (row 5) to the "Neetcode 150" sheet for the newly solved
# problem, mirroring the structure/formatting of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CATEGORY / DIFFICULTY columns - same values as the other Tree/Medium rows.
$ws.Range("A5").Value = "Tree"
$ws.Range("B5").Value = "Medium"

# NOTES column (set before NAME so the shared-string table gets the same
# ordering as the authored workbook: the notes string lands right after
# the previous row's notes, before the new problem title string).
$ws.Range("D5").Value = "Just iterate over it in dfs and keep track of the max value of all parents/ancestors of this node, if node.val >= then we can add (and update max), otherwise continue dfs"

# NAME column with a hyperlink to the Leetcode problem, matching the style
# used for the other problem-name cells (the "Good" green cell style) and
# with the display text set to the target URL like the other hyperlinks.
$url = "https://leetcode.com/problems/count-good-nodes-in-binary-tree/"
$ws.Hyperlinks.Add($ws.Range("C5"), $url, [Type]::Missing, [Type]::Missing, $url) | Out-Null
$ws.Range("C5").Value = "1448. Count Good Nodes in Binary Tree"
$ws.Range("C5").Style = "Good"

# Row height to match the other wrapped-text rows (28.8pt).
$ws.Rows.Item(5).RowHeight = 28.8

# Selection moves to D5, as in the saved workbook.
$ws.Range("D5").Select() | Out-Null
